# Driver Validation Suite - USART_Server rework
# - Remove the old "RXn/TXn" label textbox on the USART driver-validation slide
# - Add the new "USART_Server /" and "USARTn" label textboxes in its place
# - Refresh the cached "Last printed/updated" date field text on the
#   handout master and the notes master

$p = $ppt.ActivePresentation

# EMU -> point helpers (Shapes.AddTextbox / TextFrame margins are expressed
# in points; 1 pt = 12700 EMU)
function EmuToPt([double]$emu) { return $emu / 12700.0 }

# ----------------------------------------------------------------------
# 1. Slide 4: swap the "TextBox 82" (RXn/TXn) label for the two new labels
# ----------------------------------------------------------------------
$slide = $p.Slides.Item(4)

$old = $null
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.Name -eq "TextBox 82") {
        $old = $shp
        break
    }
}
if ($old -ne $null) {
    $old.Delete()
}

# --- New shape: "TextBox 146" -> "USART_Server /" ----------------------
$tb1 = $slide.Shapes.AddTextbox(
    1,
    (EmuToPt -47611),
    (EmuToPt 2641770),
    (EmuToPt 1305292),
    (EmuToPt 276999)
)
$tb1.Name = "TextBox 146"
$tb1.Fill.Visible = $false
$tb1.Line.Visible = $false
$tb1.TextFrame.WordWrap = $true
$tb1.TextFrame.AutoSize = 1
$tb1.TextFrame.MarginLeft = (EmuToPt 216000)
$tb1.TextFrame.MarginRight = (EmuToPt 36000)

$tr1 = $tb1.TextFrame.TextRange
$tr1.Text = "USART_Server /"
$tr1.Font.Size = 12
$tr1.Font.Name = "Calibri"
$tr1.Font.NameComplexScript = "Calibri"
$tr1.ParagraphFormat.Alignment = 3
$tr1.Characters(1, 13).Font.Bold = $true

# --- New shape: "TextBox 147" -> "USART" + italic "n" -------------------
$tb2 = $slide.Shapes.AddTextbox(
    1,
    (EmuToPt 312489),
    (EmuToPt 2448665),
    (EmuToPt 1011880),
    (EmuToPt 276999)
)
$tb2.Name = "TextBox 147"
$tb2.Fill.Visible = $false
$tb2.Line.Visible = $false
$tb2.TextFrame.WordWrap = $true
$tb2.TextFrame.AutoSize = 1
$tb2.TextFrame.MarginLeft = (EmuToPt 216000)
$tb2.TextFrame.MarginRight = (EmuToPt 36000)

$tr2 = $tb2.TextFrame.TextRange
$tr2.Text = "USARTn"
$tr2.Font.Size = 12
$tr2.Font.Name = "Calibri"
$tr2.Font.NameComplexScript = "Calibri"
$tr2.ParagraphFormat.Alignment = 3
$tr2.Characters(6, 1).Font.Italic = $true

# ----------------------------------------------------------------------
# 2. Refresh the cached date text on the handout master & notes master
# ----------------------------------------------------------------------
$hm = $p.HandoutMaster
$hm.HeadersFooters.DateAndTime.Text = "01-Jul-21"

$nm = $p.NotesMaster
$nm.HeadersFooters.DateAndTime.Text = "01-Jul-21"

Write-Host "done"
